# Update crypto price/volume data as scraped on Fri Nov  8 17:53:18 UTC 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "76.343.44"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.06%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.915.32"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.97%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "198.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.98%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "593.83"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.01%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("E8").Value = "  -1.23%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.193"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.86%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "2.910.23"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.79%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.447"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +14.50%  "

$ws.Range("E12").Value = "  +0.49%  "

$ws.Range("E13").Value = "  -0.48%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.448.96"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.82%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "76.103.21"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.06%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.04"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.81%  "

$ws.Range("E17").Value = "  -1.55%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.903.94"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.67%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.23"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.07%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.60"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.34%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "368.64"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.77%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.29"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.44%  "

$ws.Range("E23").Value = "  -4.09%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.11%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.11%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.045.86"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.26%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.23"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.05%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.60"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.52%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000105"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.23%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.997"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.42%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.06"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.33%  "

$ws.Range("E32").Value = "  -3.83%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "493.42"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.10%  "

$ws.Range("E34").Value = "  -0.64%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.14%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "166.16"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.24%  "

$ws.Range("E37").Value = "  -0.31%  "

$ws.Range("E38").Value = "  +11.72%  "

$ws.Range("E39").Value = "  +20.50%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.71"
$ws.Range("D40").Style = "Normal"

$ws.Range("E41").Value = "  +0.02%  "

$ws.Range("E42").Value = "  -8.24%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "178.31"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.58%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.86"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.80%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.63"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.62%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.12"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.49%  "

$ws.Range("E47").Value = "  -5.72%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.582"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.77%  "

$ws.Range("E49").Value = "  +2.01%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.26"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.97%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.43"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.85%  "
